# Update cryptocurrency price (D) and 1h volume change (E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell that carries the default (unstyled) cell format used to
# reset formatting on cells whose NumberFormat we temporarily change to text
$defaultStyle = $ws.Range("D4").Style

$updates = @(
    @{ Cell = "D2"; Value = '69.787.04'; Numeric = $false },
    @{ Cell = "E2"; Value = '  +2.92%  '; Numeric = $false },
    @{ Cell = "D3"; Value = '3.385.06'; Numeric = $false },
    @{ Cell = "E3"; Value = '  +3.76%  '; Numeric = $false },
    @{ Cell = "E4"; Value = '  -0.01%  '; Numeric = $false },
    @{ Cell = "D5"; Value = '191.51'; Numeric = $true },
    @{ Cell = "E5"; Value = '  +3.21%  '; Numeric = $false },
    @{ Cell = "D6"; Value = '594.78'; Numeric = $true },
    @{ Cell = "E6"; Value = '  +2.31%  '; Numeric = $false },
    @{ Cell = "E7"; Value = '  +0.73%  '; Numeric = $false },
    @{ Cell = "E8"; Value = '  +0.01%  '; Numeric = $false },
    @{ Cell = "E9"; Value = '  +1.93%  '; Numeric = $false },
    @{ Cell = "E10"; Value = '  +2.52%  '; Numeric = $false },
    @{ Cell = "D11"; Value = '0.420'; Numeric = $true },
    @{ Cell = "E11"; Value = '  +1.69%  '; Numeric = $false },
    @{ Cell = "D12"; Value = '3.974.45'; Numeric = $false },
    @{ Cell = "E12"; Value = '  +3.95%  '; Numeric = $false },
    @{ Cell = "E13"; Value = '  -0.67%  '; Numeric = $false },
    @{ Cell = "D14"; Value = '28.75'; Numeric = $true },
    @{ Cell = "E14"; Value = '  +3.53%  '; Numeric = $false },
    @{ Cell = "D15"; Value = '69.763.14'; Numeric = $false },
    @{ Cell = "E15"; Value = '  +2.94%  '; Numeric = $false },
    @{ Cell = "E16"; Value = '  +1.44%  '; Numeric = $false },
    @{ Cell = "D17"; Value = '3.380.05'; Numeric = $false },
    @{ Cell = "E17"; Value = '  +3.46%  '; Numeric = $false },
    @{ Cell = "D18"; Value = '455.17'; Numeric = $true },
    @{ Cell = "E18"; Value = '  +15.38%  '; Numeric = $false },
    @{ Cell = "D19"; Value = '5.83'; Numeric = $true },
    @{ Cell = "E19"; Value = '  +1.13%  '; Numeric = $false },
    @{ Cell = "D20"; Value = '13.83'; Numeric = $true },
    @{ Cell = "E20"; Value = '  +1.58%  '; Numeric = $false },
    @{ Cell = "D21"; Value = '7.81'; Numeric = $true },
    @{ Cell = "E21"; Value = '  +2.47%  '; Numeric = $false },
    @{ Cell = "D22"; Value = '76.14'; Numeric = $true },
    @{ Cell = "E22"; Value = '  +6.37%  '; Numeric = $false },
    @{ Cell = "E23"; Value = '  +0.17%  '; Numeric = $false },
    @{ Cell = "D24"; Value = '0.522'; Numeric = $true },
    @{ Cell = "E24"; Value = '  +0.79%  '; Numeric = $false },
    @{ Cell = "E25"; Value = '  +3.20%  '; Numeric = $false },
    @{ Cell = "E26"; Value = '  +2.05%  '; Numeric = $false },
    @{ Cell = "D27"; Value = '9.57'; Numeric = $true },
    @{ Cell = "E27"; Value = '  -0.58%  '; Numeric = $false },
    @{ Cell = "E28"; Value = '  -0.12%  '; Numeric = $false },
    @{ Cell = "D29"; Value = '2.03'; Numeric = $true },
    @{ Cell = "E29"; Value = '  +3.69%  '; Numeric = $false },
    @{ Cell = "D30"; Value = '23.48'; Numeric = $true },
    @{ Cell = "E30"; Value = '  +3.38%  '; Numeric = $false },
    @{ Cell = "D31"; Value = '5.63'; Numeric = $true },
    @{ Cell = "E31"; Value = '  +1.36%  '; Numeric = $false },
    @{ Cell = "E32"; Value = '  +2.54%  '; Numeric = $false },
    @{ Cell = "D33"; Value = '7.02'; Numeric = $true },
    @{ Cell = "E33"; Value = '  +0.00%  '; Numeric = $false },
    @{ Cell = "E34"; Value = '  +0.00%  '; Numeric = $false },
    @{ Cell = "E35"; Value = '  +6.06%  '; Numeric = $false },
    @{ Cell = "D36"; Value = '164.56'; Numeric = $true },
    @{ Cell = "E36"; Value = '  +1.17%  '; Numeric = $false },
    @{ Cell = "E37"; Value = '  +2.16%  '; Numeric = $false },
    @{ Cell = "D38"; Value = '27.87'; Numeric = $true },
    @{ Cell = "E38"; Value = '  +4.44%  '; Numeric = $false },
    @{ Cell = "E39"; Value = '  +0.55%  '; Numeric = $false },
    @{ Cell = "D40"; Value = '4.61'; Numeric = $true },
    @{ Cell = "E40"; Value = '  +1.29%  '; Numeric = $false },
    @{ Cell = "D41"; Value = '6.57'; Numeric = $true },
    @{ Cell = "E41"; Value = '  +1.37%  '; Numeric = $false },
    @{ Cell = "D42"; Value = '2.746.70'; Numeric = $false },
    @{ Cell = "E43"; Value = '  +1.40%  '; Numeric = $false },
    @{ Cell = "D44"; Value = '25.57'; Numeric = $true },
    @{ Cell = "E44"; Value = '  +2.94%  '; Numeric = $false },
    @{ Cell = "E45"; Value = '  -0.10%  '; Numeric = $false },
    @{ Cell = "D46"; Value = '41.16'; Numeric = $true },
    @{ Cell = "E46"; Value = '  +1.27%  '; Numeric = $false },
    @{ Cell = "D47"; Value = '339.52'; Numeric = $true },
    @{ Cell = "E47"; Value = '  +1.49%  '; Numeric = $false },
    @{ Cell = "E48"; Value = '  +2.03%  '; Numeric = $false },
    @{ Cell = "D49"; Value = '32.94'; Numeric = $true },
    @{ Cell = "E49"; Value = '  +6.46%  '; Numeric = $false },
    @{ Cell = "E50"; Value = '  +4.58%  '; Numeric = $false }
)

foreach ($item in $updates) {
    $range = $ws.Range($item.Cell)
    if ($item.Numeric) {
        # Force the cell to be stored as text so numeric-looking strings
        # (e.g. "191.51") are not reinterpreted as actual numbers.
        $range.NumberFormat = "@"
        $range.Value = $item.Value
        $range.Style = $defaultStyle
    } else {
        $range.Value = $item.Value
    }
}
